# Append two new match rows (130 and 131) to the Romania Liga-1 2023-2024
# results sheet, mirroring the existing table's layout/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (129) down
# onto the two new rows so that column A keeps its bold/bordered "index"
# style and column E keeps its date/time number format, exactly like
# every other row in the table.
$ws.Range("A129:V129").Copy()
$ws.Range("A130:V131").PasteSpecial(-4122)

# ---- Row 130 ------------------------------------------------------------
$ws.Cells.Item(130, 1).Value  = 129
$ws.Cells.Item(130, 2).Value  = "romania"
$ws.Cells.Item(130, 3).Value  = "liga-1"
$ws.Cells.Item(130, 4).Value  = "2023-2024"
$ws.Cells.Item(130, 5).Value  = 45254.6875
$ws.Cells.Item(130, 6).Value  = "FC Hermannstadt"
$ws.Cells.Item(130, 7).Value  = 0
$ws.Cells.Item(130, 8).Value  = "Poli Iasi"
$ws.Cells.Item(130, 9).Value  = 0
$ws.Cells.Item(130, 10).Value = 1.79
$ws.Cells.Item(130, 11).Value = "22/11/2023 16:42"
$ws.Cells.Item(130, 12).Value = 1.76
$ws.Cells.Item(130, 13).Value = "24/11/2023 16:29"
$ws.Cells.Item(130, 14).Value = 3.32
$ws.Cells.Item(130, 15).Value = "22/11/2023 16:42"
$ws.Cells.Item(130, 16).Value = 3.56
$ws.Cells.Item(130, 17).Value = "24/11/2023 16:29"
$ws.Cells.Item(130, 18).Value = 4.59
$ws.Cells.Item(130, 19).Value = "22/11/2023 16:42"
$ws.Cells.Item(130, 20).Value = 5.01
$ws.Cells.Item(130, 21).Value = "24/11/2023 16:29"
$ws.Cells.Item(130, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/fc-hermannstadt-poli-iasi/tzxCfe1T/"

# ---- Row 131 ------------------------------------------------------------
$ws.Cells.Item(131, 1).Value  = 130
$ws.Cells.Item(131, 2).Value  = "romania"
$ws.Cells.Item(131, 3).Value  = "liga-1"
$ws.Cells.Item(131, 4).Value  = "2023-2024"
$ws.Cells.Item(131, 5).Value  = 45254.8125
$ws.Cells.Item(131, 6).Value  = "Petrolul"
$ws.Cells.Item(131, 7).Value  = 2
$ws.Cells.Item(131, 8).Value  = "FC Botosani"
$ws.Cells.Item(131, 9).Value  = 1
$ws.Cells.Item(131, 10).Value = 1.68
$ws.Cells.Item(131, 11).Value = "23/11/2023 11:42"
$ws.Cells.Item(131, 12).Value = 1.83
$ws.Cells.Item(131, 13).Value = "24/11/2023 19:28"
$ws.Cells.Item(131, 14).Value = 3.7
$ws.Cells.Item(131, 15).Value = "23/11/2023 11:42"
$ws.Cells.Item(131, 16).Value = 3.44
$ws.Cells.Item(131, 17).Value = "24/11/2023 19:28"
$ws.Cells.Item(131, 18).Value = 5.19
$ws.Cells.Item(131, 19).Value = "23/11/2023 11:42"
$ws.Cells.Item(131, 20).Value = 4.75
$ws.Cells.Item(131, 21).Value = "24/11/2023 19:28"
$ws.Cells.Item(131, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/petrolul-fc-botosani/ltawBGoA/"
